$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need to be forced to
# text (matching the source inline-string cell type) since Excel would
# otherwise auto-convert them to numeric values on assignment.
$textForceCells = @(
    'D4',
    'D5',
    'D7',
    'D8',
    'D9',
    'D10',
    'D11',
    'D13',
    'D14',
    'D15',
    'D16',
    'D17',
    'D18',
    'D22',
    'D23',
    'D24',
    'D25',
    'D26',
    'D27',
    'D28',
    'D29',
    'D30',
    'D33',
    'D35',
    'D36',
    'D37',
    'D38',
    'D39',
    'D40',
    'D41',
    'D42',
    'D43',
    'D44',
    'D45',
    'D46',
    'D47',
    'D48',
    'D50',
    'D51',
)
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated Price (D) / Volume(1h) (E) values cell by cell.
$ws.Range('D2').Value = '27.600.94'
$ws.Range('E2').Value = '  +2.45%  '
$ws.Range('D3').Value = '1.854.64'
$ws.Range('E3').Value = '  +2.35%  '
$ws.Range('D4').Value = '1.033'
$ws.Range('E4').Value = '  +2.78%  '
$ws.Range('D5').Value = '322.42'
$ws.Range('E5').Value = '  +3.54%  '
$ws.Range('E6').Value = '  +2.44%  '
$ws.Range('D7').Value = '0.4396'
$ws.Range('E7').Value = '  +2.43%  '
$ws.Range('D8').Value = '0.3788'
$ws.Range('E8').Value = '  +2.61%  '
$ws.Range('D9').Value = '0.07416'
$ws.Range('E9').Value = '  +2.29%  '
$ws.Range('D10').Value = '0.8791'
$ws.Range('E10').Value = '  +1.97%  '
$ws.Range('D11').Value = '21.65'
$ws.Range('E11').Value = '  +2.38%  '
$ws.Range('D12').Value = '1.868.24'
$ws.Range('E12').Value = '  -7.47%  '
$ws.Range('D13').Value = '5.533'
$ws.Range('E13').Value = '  +2.62%  '
$ws.Range('D14').Value = '6.710'
$ws.Range('E14').Value = '  +1.22%  '
$ws.Range('D15').Value = '0.07213'
$ws.Range('E15').Value = '  +4.65%  '
$ws.Range('D16').Value = '83.31'
$ws.Range('E16').Value = '  +3.17%  '
$ws.Range('D17').Value = '1.035'
$ws.Range('E17').Value = '  +2.99%  '
$ws.Range('D18').Value = '0.000009072'
$ws.Range('E18').Value = '  +2.48%  '
$ws.Range('E19').Value = '  +2.46%  '
$ws.Range('E20').Value = '  +1.83%  '
$ws.Range('D21').Value = '27.626.98'
$ws.Range('E21').Value = '  +2.49%  '
$ws.Range('D22').Value = '5.282'
$ws.Range('E22').Value = '  +1.81%  '
$ws.Range('D23').Value = '11.42'
$ws.Range('E23').Value = '  +3.96%  '
$ws.Range('D24').Value = '157.98'
$ws.Range('E24').Value = '  +2.76%  '
$ws.Range('D25').Value = '1.916'
$ws.Range('E25').Value = '  +1.73%  '
$ws.Range('D26').Value = '18.76'
$ws.Range('E26').Value = '  +2.73%  '
$ws.Range('D27').Value = '1.980'
$ws.Range('E27').Value = '  +4.16%  '
$ws.Range('D28').Value = '5.298'
$ws.Range('E28').Value = '  +1.47%  '
$ws.Range('D29').Value = '117.44'
$ws.Range('E29').Value = '  +2.10%  '
$ws.Range('D30').Value = '0.09064'
$ws.Range('E30').Value = '  +1.49%  '
$ws.Range('E31').Value = '  +3.85%  '
$ws.Range('E32').Value = '  +2.72%  '
$ws.Range('D33').Value = '4.544'
$ws.Range('E33').Value = '  +2.73%  '
$ws.Range('E34').Value = '  +3.00%  '
$ws.Range('D35').Value = '1.029'
$ws.Range('E35').Value = '  +1.99%  '
$ws.Range('D36').Value = '1.155'
$ws.Range('E36').Value = '  +2.90%  '
$ws.Range('D37').Value = '0.01977'
$ws.Range('E37').Value = '  +3.01%  '
$ws.Range('D38').Value = '0.05315'
$ws.Range('E38').Value = '  +2.05%  '
$ws.Range('D39').Value = '0.5169'
$ws.Range('E39').Value = '  +1.75%  '
$ws.Range('D40').Value = '2.822'
$ws.Range('E40').Value = '  +3.33%  '
$ws.Range('D41').Value = '0.1683'
$ws.Range('E41').Value = '  +2.25%  '
$ws.Range('D42').Value = '6.766'
$ws.Range('E42').Value = '  +5.34%  '
$ws.Range('D43').Value = '8.601'
$ws.Range('E43').Value = '  +3.89%  '
$ws.Range('D44').Value = '109.01'
$ws.Range('E44').Value = '  +1.93%  '
$ws.Range('D45').Value = '10.56'
$ws.Range('E45').Value = '  +1.92%  '
$ws.Range('D46').Value = '1.715'
$ws.Range('E46').Value = '  +3.86%  '
$ws.Range('D47').Value = '0.4663'
$ws.Range('E47').Value = '  +1.78%  '
$ws.Range('D48').Value = '0.06407'
$ws.Range('E48').Value = '  +2.14%  '
$ws.Range('E49').Value = '  +2.70%  '
$ws.Range('D50').Value = '39.49'
$ws.Range('E50').Value = '  +4.69%  '
$ws.Range('D51').Value = '64.15'
$ws.Range('E51').Value = '  +0.86%  '

# Restore the default (unstyled) cell style on the cells we text-forced,
# so only the value changes -- not the formatting/style of the cell.
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).Style = "Normal"
}
